# Update the "Admin" test-data worksheet with rotated test credentials:
#   - D3 (hyperlinked site URL): trim the "/ioadmin/" admin path off the URL
#   - B8 (login user):   manisha@vtestcorp.com  -> saumyata@vtestcorp.com
#   - C8 (login pwd):    Manisha1!              -> Saumyata2!
# Also add a hyperlink on D3 pointing at the (now-shortened) URL, matching
# the existing hyperlink already present on B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- update cell values -------------------------------------------------
$ws.Range("D3").Value = "https://www.primebpmapp.com/"
$ws.Range("B8").Value = "saumyata@vtestcorp.com"
$ws.Range("C8").Value = "Saumyata2!"

# --- add hyperlink on D3 --------------------------------------------------
# Pre-apply the built-in "Hyperlink" style (D3 already renders with the
# hyperlink font) so that adding the hyperlink does not fabricate a brand
# new cell style entry in styles.xml.
$ws.Range("D3").Style = "Hyperlink"
[void]$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.primebpmapp.com/")

# --- restore the saved cursor/selection state ----------------------------
[void]$ws.Activate()
[void]$ws.Range("D19").Select()
